$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same shared style (bold, centered,
# bordered) used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels for the new "season record" columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-42) gets the same team season record values.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
